$d = $word.ActiveDocument

$replacements = @(
    ,@('48÷5=9, 3', '66÷9=7, 3')
    ,@('96÷9=10, 6', '55÷5=11, 0')
    ,@('49÷5=9, 4', '22÷6=3, 4')
    ,@('52÷2=26, 0', '53÷6=8, 5')
    ,@('29÷9=3, 2', '41÷8=5, 1')
    ,@('27÷9=3, 0', '82÷3=27, 1')
    ,@('77÷8=9, 5', '39÷6=6, 3')
    ,@('70÷3=23, 1', '18÷5=3, 3')
    ,@('49÷3=16, 1', '80÷8=10, 0')
    ,@('78÷8=9, 6', '67÷5=13, 2')
    ,@('16÷5=3, 1', '78÷9=8, 6')
    ,@('25÷2=12, 1', '10÷2=5, 0')
    ,@('38÷2=19, 0', '36÷6=6, 0')
    ,@('77÷9=8, 5', '23÷4=5, 3')
    ,@('69÷8=8, 5', '63÷3=21, 0')
    ,@('24÷2=12, 0', '68÷6=11, 2')
    ,@('41÷4=10, 1', '63÷4=15, 3')
    ,@('11÷4=2, 3', '62÷8=7, 6')
    ,@('14÷9=1, 5', '11÷9=1, 2')
    ,@('92÷4=23, 0', '51÷4=12, 3')
    ,@('77÷5=15, 2', '82÷2=41, 0')
    ,@('28÷5=5, 3', '90÷7=12, 6')
    ,@('17÷6=2, 5', '12÷6=2, 0')
    ,@('66÷5=13, 1', '90÷9=10, 0')
    ,@('43÷5=8, 3', '54÷3=18, 0')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find text: $old"
    }
}

$d.Save()
